# Fix the #REF! formulas in row 6 of the "save to libraries.csv" sheet so
# that they once again point at the "Batik" row (row 7) of the "libraries"
# sheet, and update the active sheet / selection to match.

$wb = $excel.ActiveWorkbook

$wsCsv = $wb.Worksheets.Item("save to libraries.csv")

# Repair the broken #REF! formulas in row 6 (A6:D6) of the csv sheet so
# they reference the corresponding cells of row 7 on the libraries sheet.
$wsCsv.Range("A6").Formula = "=libraries!A7"
$wsCsv.Range("B6").Formula = "=libraries!B7"
$wsCsv.Range("C6").Formula = "=libraries!C7"
$wsCsv.Range("D6").Formula = "=libraries!D7"

# Make "save to libraries.csv" the active sheet, with A20:D20 selected.
$wsCsv.Select()
$wsCsv.Range("A20:D20").Select()
